$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "41.567.93"
Set-TextValue "E2" "  +0.09%  "
Set-TextValue "D3" "2.455.62"
Set-TextValue "E3" "  -1.61%  "
Set-TextValue "E4" "  +0.37%  "
Set-TextValue "D5" "314.50"
Set-TextValue "E5" "  +1.06%  "
Set-TextValue "D6" "92.17"
Set-TextValue "E6" "  -0.63%  "
Set-TextValue "E7" "  +1.28%  "
Set-TextValue "E8" "  +0.19%  "
Set-TextValue "D9" "0.508"
Set-TextValue "E9" "  +2.80%  "
Set-TextValue "D10" "32.27"
Set-TextValue "E10" "  -0.08%  "
Set-TextValue "D11" "0.0795"
Set-TextValue "E11" "  +2.16%  "
Set-TextValue "E12" "  +0.66%  "
Set-TextValue "D13" "2.835.02"
Set-TextValue "D14" "6.82"
Set-TextValue "E14" "  -0.02%  "
Set-TextValue "E15" "  +2.43%  "
Set-TextValue "D16" "2.483.99"
Set-TextValue "E16" "  -0.08%  "
Set-TextValue "D17" "0.773"
Set-TextValue "E17" "  +1.35%  "
Set-TextValue "D18" "41.559.68"
Set-TextValue "E18" "  +0.14%  "
Set-TextValue "E19" "  +2.12%  "
Set-TextValue "E20" "  +1.25%  "
Set-TextValue "D21" "70.73"
Set-TextValue "E21" "  +0.23%  "
Set-TextValue "D22" "11.31"
Set-TextValue "E22" "  +1.49%  "
Set-TextValue "D23" "237.86"
Set-TextValue "E23" "  +1.15%  "
Set-TextValue "D24" "2.69"
Set-TextValue "E24" "  -0.36%  "
Set-TextValue "E25" "  -0.07%  "
Set-TextValue "D26" "1.90"
Set-TextValue "E26" "  +0.19%  "
Set-TextValue "D27" "24.27"
Set-TextValue "E27" "  -0.60%  "
Set-TextValue "D28" "2.25"
Set-TextValue "E28" "  +0.52%  "
Set-TextValue "D29" "9.65"
Set-TextValue "E29" "  +0.24%  "
Set-TextValue "D30" "34.92"
Set-TextValue "E30" "  -4.38%  "
Set-TextValue "D31" "155.71"
Set-TextValue "E31" "  +1.46%  "
Set-TextValue "E32" "  +0.85%  "
Set-TextValue "D33" "2.56"
Set-TextValue "E33" "  +0.26%  "
Set-TextValue "D34" "0.0758"
Set-TextValue "E34" "  -0.21%  "
Set-TextValue "D35" "2.49"
Set-TextValue "E35" "  +0.07%  "
Set-TextValue "D36" "17.40"
Set-TextValue "E36" "  -4.66%  "
Set-TextValue "E37" "  -3.88%  "
Set-TextValue "E38" "  +1.23%  "
Set-TextValue "E39" "  +0.55%  "
Set-TextValue "D40" "1.79"
Set-TextValue "E40" "  -3.09%  "
Set-TextValue "D41" "3.94"
Set-TextValue "E41" "  -3.83%  "
Set-TextValue "E42" "  -0.23%  "
Set-TextValue "D43" "1.966.08"
Set-TextValue "E43" "  +0.89%  "
Set-TextValue "E46" "  -2.28%  "
Set-TextValue "E47" "  +2.05%  "
Set-TextValue "D48" "2.694.40"
Set-TextValue "E48" "  -1.05%  "
Set-TextValue "D49" "96.27"
Set-TextValue "E49" "  +0.34%  "
Set-TextValue "D50" "66.37"
Set-TextValue "E50" "  -0.99%  "
Set-TextValue "E51" "  -2.81%  "

# Row 44/45 swap (VeChain/EnergySwap reorder with updated values)
Set-TextValue "B44" "VeChain"
Set-TextValue "C44" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D44" "0.0281"
Set-TextValue "E44" "  -0.11%  "
Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "18.64"
Set-TextValue "E45" "  -5.48%  "
